$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.798.06"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "2.293.51"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'116.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +12.68%  "

$ws.Range("D6").Value = "'268.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("E9").Value = "  +1.54%  "

$ws.Range("D10").Value = "'49.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.47%  "

$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("D12").Value = "'8.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +12.11%  "

$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("D14").Value = "'15.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.85%  "

$ws.Range("D15").Value = "2.638.69"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("D16").Value = "'0.882"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.97%  "

$ws.Range("D17").Value = "2.277.41"
$ws.Range("E17").Value = "  -1.10%  "

$ws.Range("D18").Value = "43.726.20"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("E19").Value = "  -1.23%  "

$ws.Range("E20").Value = "  +11.85%  "

$ws.Range("D21").Value = "'72.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "'2.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.52%  "

$ws.Range("D23").Value = "'9.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.13%  "

$ws.Range("D24").Value = "'233.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("E25").Value = "  +2.54%  "

$ws.Range("D27").Value = "'11.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.11%  "

$ws.Range("E28").Value = "  +1.99%  "

$ws.Range("D29").Value = "'42.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.85%  "

$ws.Range("D30").Value = "'3.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.87%  "

$ws.Range("E31").Value = "  -2.25%  "

$ws.Range("D32").Value = "'173.27"
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").Value = "'0.0934"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.77%  "

$ws.Range("D34").Value = "'21.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.86%  "

$ws.Range("E35").Value = "  +5.18%  "

$ws.Range("E36").Value = "  +0.29%  "

$ws.Range("D37").Value = "'4.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.04%  "

$ws.Range("E38").Value = "  +1.65%  "

$ws.Range("D39").Value = "'0.108"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("E40").Value = "  +6.93%  "

$ws.Range("D41").Value = "'14.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.83%  "

$ws.Range("E42").Value = "  +3.90%  "

$ws.Range("D43").Value = "'74.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.71%  "

$ws.Range("E44").Value = "  +2.74%  "

$ws.Range("D45").Value = "'6.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +22.39%  "

$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").Value = "'1.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.30%  "

$ws.Range("D48").Value = "'8.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("E49").Value = "  +4.46%  "

$ws.Range("D50").Value = "'102.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.54%  "

$ws.Range("E51").Value = "  -1.46%  "
